$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell L1: "optimierungen" (percentage-style header, like F1:J1) ---
$ws.Range("L1").Value = "optimierungen"
$ws.Range("L1").NumberFormat = "0%"

# --- Row 4 gets a new "put+take" tag in the new column ---
$ws.Range("L4").Value = "put+take"

# --- Fill in the previously-incomplete row 5 with the executed "chw-work" experiment ---
$ws.Range("A5").Value = "chw-work"
$ws.Range("B5").Value = "analysis.throughput (timestamp)"
$ws.Range("C5").Value = 100000
$ws.Range("D5").Value = 800
$ws.Range("E5").Value = 634
$ws.Range("F5").Value = 453
$ws.Range("G5").Value = 475
$ws.Range("H5").Value = 493
$ws.Range("I5").Value = 564
$ws.Range("J5").Value = 3237
$ws.Range("K5").Value = 2.914
$ws.Range("L5").Value = "put+take"

# --- Grow Table2 to include the new "optimierungen" column ---
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Add()
$col.Range.Cells.Item(1, 1).Value = "optimierungen"

# --- Move the active selection to reflect where work continued ---
$ws.Range("B14").Select()
